$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "07/08/2025"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "Port Vale"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "Cardiff"
$ws.Range("F3").Value = "D"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.89
$ws.Range("L3").Value = 1.5
$ws.Range("M3").Value = 9
$ws.Range("N3").Value = 16
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
